$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 10442.454
$ws.Range("I6").Value = 1207.5555
$ws.Range("K6").Value = 3622.6665
$ws.Range("M6").Value = -3510.6665

$ws.Range("H11").Value = 61113744
$ws.Range("I11").Value = 61113744
$ws.Range("K11").Value = 61113744
$ws.Range("M11").Value = -61113604

$ws.Range("H107").Value = 10882.9
$ws.Range("I107").Value = 11992.111
$ws.Range("K107").Value = 11992.111
$ws.Range("M107").Value = -10072.111

$ws.Range("H113").Value = 2489.1892
$ws.Range("I113").Value = 2185.7144
$ws.Range("J113").Value = 2673.913
$ws.Range("K113").Value = 2185.7144
$ws.Range("L113").Value = 2673.913
$ws.Range("M113").Value = 1068.2856
$ws.Range("N113").Value = -9181.913

$ws.Range("H117").Value = 33230
$ws.Range("J117").Value = 33230
$ws.Range("L117").Value = 33230
$ws.Range("N117").Value = -42408

$ws.Range("H131").Value = 2591.4614
$ws.Range("J131").Value = 4875
$ws.Range("L131").Value = 14625
$ws.Range("N131").Value = -24705

$ws.Range("H132").Value = 35185.1
$ws.Range("I132").Value = 5903.421
$ws.Range("J132").Value = 90820.3
$ws.Range("K132").Value = 17710.263
$ws.Range("L132").Value = 272460.9
$ws.Range("M132").Value = -15180.263
$ws.Range("N132").Value = -277520.9

$ws.Range("H139").Value = 32342.5
$ws.Range("J139").Value = 32342.5
$ws.Range("L139").Value = 32342.5
$ws.Range("N139").Value = -42622.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 83826.5
$ws.Range("J103").Value = 83826.5
$ws.Range("L103").Value = 83826.5
$ws.Range("N103").Value = -86170.5

$ws.Range("H111").Value = 47698
$ws.Range("J111").Value = 47698
$ws.Range("L111").Value = 47698
$ws.Range("N111").Value = -55878

$ws.Range("H135").Value = 25154.545
$ws.Range("J135").Value = 25154.545
$ws.Range("L135").Value = 25154.545
$ws.Range("N135").Value = -35294.545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2669.3333
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 2008
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 2008
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -4254

$ws.Range("H89").Value = 2669.3333
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 2008
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 10040
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -21272

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H115").Value = 28958.334
$ws.Range("J115").Value = 28958.334
$ws.Range("L115").Value = 28958.334
$ws.Range("N115").Value = -31308.334

$ws.Range("H116").Value = 43134.8
$ws.Range("J116").Value = 43134.8
$ws.Range("L116").Value = 43134.8
$ws.Range("N116").Value = -52312.8

$ws.Range("H118").Value = 44716
$ws.Range("J118").Value = 44716
$ws.Range("L118").Value = 44716
$ws.Range("N118").Value = -48030

$ws.Range("H120").Value = 32613.727
$ws.Range("J120").Value = 32613.727
$ws.Range("L120").Value = 32613.727
$ws.Range("N120").Value = -39871.727

$ws.Range("H133").Value = 22437.6
$ws.Range("J133").Value = 22437.6
$ws.Range("L133").Value = 22437.6
$ws.Range("N133").Value = -27497.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1453
$ws.Range("I7").Value = 3533
$ws.Range("J7").Value = 413
$ws.Range("K7").Value = 10599
$ws.Range("L7").Value = 1239
$ws.Range("M7").Value = -10487
$ws.Range("N7").Value = -1463

$ws.Range("H92").Value = 831.9545000000001
$ws.Range("I92").Value = 676.1177
$ws.Range("J92").Value = 1361.8
$ws.Range("K92").Value = 2028.3531
$ws.Range("L92").Value = 4085.4
$ws.Range("M92").Value = -780.3531
$ws.Range("N92").Value = -6581.4

$ws.Range("H120").Value = 378557.38
$ws.Range("I120").Value = 501076.66
$ws.Range("J120").Value = 10999.5
$ws.Range("K120").Value = 1503229.98
$ws.Range("L120").Value = 32998.5
$ws.Range("M120").Value = -1498391.98
$ws.Range("N120").Value = -42674.5

$ws.Range("H131").Value = 5948.636
$ws.Range("I131").Value = 8067.154
$ws.Range("J131").Value = 2888.5557
$ws.Range("K131").Value = 24201.462
$ws.Range("L131").Value = 8665.667099999999
$ws.Range("M131").Value = -19161.462
$ws.Range("N131").Value = -18745.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H130").Value = 45907.5
$ws.Range("J130").Value = 45907.5
$ws.Range("L130").Value = 45907.5
$ws.Range("N130").Value = -55947.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 10666.667
$ws.Range("I45").Value = 10666.667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 10666.667
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -10259.667
$ws.Range("N45").ClearContents()

$ws.Range("H76").Value = 15821.2
$ws.Range("J76").Value = 19948.334
$ws.Range("L76").Value = 19948.334
$ws.Range("N76").Value = -20624.334

$ws.Range("H79").Value = 15821.2
$ws.Range("J79").Value = 19948.334
$ws.Range("L79").Value = 19948.334
$ws.Range("N79").Value = -22288.334

$ws.Range("H110").Value = 28018
$ws.Range("J110").Value = 28018
$ws.Range("L110").Value = 28018
$ws.Range("N110").Value = -36198

$ws.Range("H111").Value = 43843
$ws.Range("J111").Value = 43843
$ws.Range("L111").Value = 43843
$ws.Range("N111").Value = -52023

$ws.Range("H112").Value = 37113.25
$ws.Range("J112").Value = 37113.25
$ws.Range("L112").Value = 37113.25
$ws.Range("N112").Value = -40067.25

$ws.Range("H116").Value = 45664
$ws.Range("J116").Value = 45664
$ws.Range("L116").Value = 45664
$ws.Range("N116").Value = -54842

$ws.Range("H118").Value = 39026.715
$ws.Range("J118").Value = 39026.715
$ws.Range("L118").Value = 39026.715
$ws.Range("N118").Value = -42340.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1421.1111
$ws.Range("I81").Value = 1441.4286
$ws.Range("J81").Value = 1350
$ws.Range("K81").Value = 2882.8572
$ws.Range("L81").Value = 2700
$ws.Range("M81").Value = -1821.8572
$ws.Range("N81").Value = -4822

$ws.Range("H84").Value = 1421.1111
$ws.Range("I84").Value = 1441.4286
$ws.Range("J84").Value = 1350
$ws.Range("K84").Value = 14414.286
$ws.Range("L84").Value = 13500
$ws.Range("M84").Value = -9110.286
$ws.Range("N84").Value = -24108

$ws.Range("H119").Value = 48694
$ws.Range("J119").Value = 48694
$ws.Range("L119").Value = 48694
$ws.Range("N119").Value = -58370

$ws.Range("H126").Value = 1280838.9
$ws.Range("I126").Value = 1472389.5
$ws.Range("K126").Value = 4417168.5
$ws.Range("M126").Value = -4414698.5

$ws.Range("H129").Value = 39425
$ws.Range("J129").Value = 39425
$ws.Range("L129").Value = 39425
$ws.Range("N129").Value = -49425
